# Price Tracker update: insert two new columns at the front (SKU Name, last-update timestamp)
# and append 25 summary rows (one per unique SKU) with Product name + latest Price.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column A; this shifts the existing
# Timestamp/Product/Price/URL/Price Change columns from A:E to C:G.
$ws.Columns("A:B").Insert()

# Copy the existing header style (bold, centered, bordered) onto the two new header cells.
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header labels.
$ws.Range("A1").Value = "SKU Name"
$ws.Range("B1").Value = "2025-12-19 06:11"

$newRows = @(
    @{Product='Jr. Sr. Baby Diaper Pants | XXL Size (15-25 kg), 42 Count | Pack of 1 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=499},
    @{Product='Jr. Sr. Baby Diaper Pants | 10-12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy (XL, 28, count)'; Price=299},
    @{Product='Jr. Sr. Baby Diaper Pants | 10-12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy (M, 36, count)'; Price=299},
    @{Product='Jr. Sr. Baby Tape Style Diaper | Adjustable Fit For Babies | Extra Large (XL) Size (12+ kg), 38 Count | Pack of 1 | Upto 12 hr Absorption with Advanced Leak Protection | Soft & Comfy'; Price=569},
    @{Product='Jr. Sr. Baby Tape Style Diaper | Adjustable Fit For Babies | Newborn/Extra Small (NB/XS) Size (0-5kg), 72 Count | Pack of 1 | Upto 12 hr Absorption with Advanced Leak Protection | Soft & Comfy'; Price=569},
    @{Product='Jr. Sr. Baby Diaper Pants | Extra Large (XL) Size (12-17 kg), 162 Count | Pack of 3 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=1299},
    @{Product='Jr. Sr. Baby Diaper Pants | Medium (M) Size (7-12 kg), 216 Count | Pack of 3 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=1299},
    @{Product='Jr. Sr. Baby Diaper Pants | XXL Size (15-25 kg), 126 Count | Pack of 3 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=1299},
    @{Product='Jr. Sr. Baby Diaper Pants | XXL Size (15-25 kg), 84 Count | Pack of 2 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=929},
    @{Product='Jr. Sr. Baby Diaper Pants | Small (S) Size (4-8 kg), 78 Count | Pack of 1 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=499},
    @{Product='Jr. Sr. Baby Tape Style Diaper | Adjustable Fit For Babies | Medium (M) Size (6-11 kg), 56 Count | Pack of 1 | Upto 12 hr Absorption with Advanced Leak Protection | Soft & Comfy'; Price=569},
    @{Product='Jr. Sr. Baby Diaper Pants | 10-12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy (XXL, 22, count)'; Price=299},
    @{Product='Jr. Sr. Baby Diaper Pants | Large (L) Size (9-14 kg), 62 Count | Pack of 1 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=499},
    @{Product='Jr. Sr. Baby Diaper Pants | 10-12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy (S, 40, count)'; Price=299},
    @{Product='Jr. Sr. Baby Diaper Pants | Small (S) Size (4-8 kg), 234 Count | Pack of 3 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=2997},
    @{Product='Jr. Sr. Baby Diaper Pants | Extra Large (XL) Size (12-17 kg), 54 Count | Pack of 1 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=499},
    @{Product='Jr. Sr. Baby Diaper Pants | Medium (M) Size (7-12 kg), 72 Count | Pack of 1 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=499},
    @{Product='Jr. Sr. Baby Diaper Pants | Medium (M) Size (7-12 kg), 144 Count | Pack of 2 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=929},
    @{Product='Jr. Sr. Baby Diaper Pants | Large (L) Size (9-14 kg), 124 Count | Pack of 2 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=929},
    @{Product='Jr. Sr. Baby Diaper Pants | 10-12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy (L, 32, count)'; Price=299},
    @{Product='Jr. Sr. Baby Diaper Pants | Large (L) Size (9-14 kg), 186 Count | Pack of 3 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=929},
    @{Product='Jr. Sr. Baby Tape Style Diaper | Adjustable Fit For Babies | Small (S) Size (3-8 kg), 64 Count | Pack of 1 | Upto 12 hr Absorption with Advanced Leak Protection | Soft & Comfy'; Price=569},
    @{Product='Jr. Sr. Baby Diaper Pants | Small (S) Size (4-8 kg), 156 Count | Pack of 2 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=929},
    @{Product='Jr. Sr. Baby Tape Style Diaper | Adjustable Fit For Babies | Large (L) Size (9-14 kg), 48 Count | Pack of 1 | Upto 12 hr Absorption with Advanced Leak Protection | Soft & Comfy'; Price=569},
    @{Product='Jr. Sr. Baby Diaper Pants | Extra Large (XL) Size (12-17 kg), 108 Count | Pack of 2 | Upto 12 hr Absorption with Advanced Leak Protection | Safe For Babies Skin | Soft, Secure & Comfy'; Price=929},
)

# Append the 25 SKU summary rows starting at row 52.
$startRow = 52
$i = 0
foreach ($item in $newRows) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $item.Product
    $ws.Cells.Item($r, 2).Value = $item.Price
    $i = $i + 1
}
